$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(125, 13).Value = -137679  # M125
$ws.Cells.Item(125, 11).Value = 140139  # K125
$ws.Cells.Item(125, 8).Value = 25844.924  # H125
$ws.Cells.Item(125, 9).Value = 15571  # I125
$ws.Cells.Item(132, 9).Value = 1541.7715  # I132
$ws.Cells.Item(132, 11).Value = 4625.3145  # K132
$ws.Cells.Item(132, 8).Value = 2274226.5  # H132
$ws.Cells.Item(132, 13).Value = -2095.3145  # M132
$ws.Cells.Item(137, 14).Value = -3634065.3  # N137
$ws.Cells.Item(137, 13).Value = -1353.2142  # M137
$ws.Cells.Item(137, 8).Value = 681000.25  # H137
$ws.Cells.Item(137, 11).Value = 3903.2142  # K137
$ws.Cells.Item(137, 10).Value = 1209655.1  # J137
$ws.Cells.Item(137, 9).Value = 1301.0714  # I137
$ws.Cells.Item(137, 12).Value = 3628965.3  # L137
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(29, 8).Value = 2500  # H29
$ws.Cells.Item(29, 14).Value = -3116  # N29
$ws.Cells.Item(29, 10).Value = 2500  # J29
$ws.Cells.Item(29, 12).Value = 2500  # L29
$ws.Cells.Item(32, 8).Value = 82051.60000000001  # H32
$ws.Cells.Item(32, 9).Value = 83177.41  # I32
$ws.Cells.Item(32, 11).Value = 83177.41  # K32
$ws.Cells.Item(32, 13).Value = -82890.41  # M32
$ws.Cells.Item(32, 12).Value = 10000  # L32
$ws.Cells.Item(32, 14).Value = -10574  # N32
$ws.Cells.Item(32, 10).Value = 10000  # J32
$ws.Cells.Item(61, 13).Value = -2903.75  # M61
$ws.Cells.Item(61, 11).Value = 3115.75  # K61
$ws.Cells.Item(61, 8).Value = 370703  # H61
$ws.Cells.Item(61, 9).Value = 3115.75  # I61
$ws.Cells.Item(122, 13).Value = -10101521.5  # M122
$ws.Cells.Item(122, 9).Value = 3367990.5  # I122
$ws.Cells.Item(122, 11).Value = 10103971.5  # K122
$ws.Cells.Item(122, 8).Value = 2416851.2  # H122
$ws.Cells.Item(136, 13).Value = -6797.25  # M136
$ws.Cells.Item(136, 9).Value = 3115.75  # I136
$ws.Cells.Item(136, 8).Value = 370703  # H136
$ws.Cells.Item(136, 11).Value = 9347.25  # K136
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(94, 13).Value = -1751.5454  # M94
$ws.Cells.Item(94, 11).Value = 2202.5454  # K94
$ws.Cells.Item(94, 9).Value = 2202.5454  # I94
$ws.Cells.Item(94, 8).Value = 2202.5454  # H94
$ws.Cells.Item(107, 13).Value = -1684.4736  # M107
$ws.Cells.Item(107, 9).Value = 3604.4736  # I107
$ws.Cells.Item(107, 11).Value = 3604.4736  # K107
$ws.Cells.Item(107, 10).Value = 6285.643  # J107
$ws.Cells.Item(107, 12).Value = 6285.643  # L107
$ws.Cells.Item(107, 8).Value = 4741.9395  # H107
$ws.Cells.Item(107, 14).Value = -10125.643  # N107
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(7, 8).Value = 161.3  # H7
$ws.Cells.Item(7, 9).Value = 201.71428  # I7
$ws.Cells.Item(7, 13).Value = -88.71428  # M7
$ws.Cells.Item(7, 11).Value = 201.71428  # K7
$ws.Cells.Item(11, 14).Value = $null  # N11
$ws.Cells.Item(11, 8).Value = 0  # H11
$ws.Cells.Item(11, 12).Value = 0  # L11
$ws.Cells.Item(11, 10).Value = 0  # J11
$ws.Cells.Item(22, 12).Value = 1948.8334  # L22
$ws.Cells.Item(22, 11).Value = 317  # K22
$ws.Cells.Item(22, 13).Value = 33  # M22
$ws.Cells.Item(22, 8).Value = 1296.1  # H22
$ws.Cells.Item(22, 9).Value = 317  # I22
$ws.Cells.Item(22, 10).Value = 1948.8334  # J22
$ws.Cells.Item(22, 14).Value = -2648.8334  # N22
$ws.Cells.Item(44, 13).Value = $null  # M44
$ws.Cells.Item(44, 8).Value = 0  # H44
$ws.Cells.Item(44, 11).Value = 0  # K44
$ws.Cells.Item(44, 12).Value = 0  # L44
$ws.Cells.Item(44, 9).Value = 0  # I44
$ws.Cells.Item(44, 10).Value = 0  # J44
$ws.Cells.Item(44, 14).Value = $null  # N44
$ws.Cells.Item(94, 13).Value = -1299  # M94
$ws.Cells.Item(94, 11).Value = 1750  # K94
$ws.Cells.Item(94, 9).Value = 1750  # I94
$ws.Cells.Item(94, 8).Value = 1752.9166  # H94
$ws.Cells.Item(132, 12).Value = 39000.429  # L132
$ws.Cells.Item(132, 9).Value = 1546.2632  # I132
$ws.Cells.Item(132, 11).Value = 4638.7896  # K132
$ws.Cells.Item(132, 10).Value = 13000.143  # J132
$ws.Cells.Item(132, 8).Value = 4630  # H132
$ws.Cells.Item(132, 13).Value = -2108.7896  # M132
$ws.Cells.Item(132, 14).Value = -44060.429  # N132
$ws.Cells.Item(141, 10).Value = 204883.36  # J141
$ws.Cells.Item(141, 8).Value = 204883.36  # H141
$ws.Cells.Item(141, 14).Value = -215243.36  # N141
$ws.Cells.Item(141, 12).Value = 204883.36  # L141
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(5, 10).Value = 1979  # J5
$ws.Cells.Item(5, 13).Value = -2671.8572  # M5
$ws.Cells.Item(5, 14).Value = -6161  # N5
$ws.Cells.Item(5, 9).Value = 927.9524  # I5
$ws.Cells.Item(5, 8).Value = 1130.0769  # H5
$ws.Cells.Item(5, 12).Value = 5937  # L5
$ws.Cells.Item(5, 11).Value = 2783.8572  # K5
$ws.Cells.Item(94, 10).Value = 12360  # J94
$ws.Cells.Item(94, 14).Value = -38432  # N94
$ws.Cells.Item(94, 8).Value = 9650  # H94
$ws.Cells.Item(94, 12).Value = 37080  # L94
$ws.Cells.Item(109, 8).Value = 11291.667  # H109
$ws.Cells.Item(109, 11).Value = 4650  # K109
$ws.Cells.Item(109, 13).Value = -3610  # M109
$ws.Cells.Item(109, 9).Value = 1550  # I109
$ws.Cells.Item(135, 14).Value = -22881  # N135
$ws.Cells.Item(135, 9).Value = 927.9524  # I135
$ws.Cells.Item(135, 8).Value = 1130.0769  # H135
$ws.Cells.Item(135, 12).Value = 17811  # L135
$ws.Cells.Item(135, 11).Value = 8351.571599999999  # K135
$ws.Cells.Item(135, 13).Value = -5816.571599999999  # M135
$ws.Cells.Item(135, 10).Value = 1979  # J135
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(107, 13).Value = 1647.25  # M107
$ws.Cells.Item(107, 9).Value = 272.75  # I107
$ws.Cells.Item(107, 11).Value = 272.75  # K107
$ws.Cells.Item(107, 10).Value = 3687.5  # J107
$ws.Cells.Item(107, 12).Value = 3687.5  # L107
$ws.Cells.Item(107, 8).Value = 2549.25  # H107
$ws.Cells.Item(107, 14).Value = -7527.5  # N107
$ws.Cells.Item(122, 14).Value = -23274.625  # N122
$ws.Cells.Item(122, 13).Value = -30814  # M122
$ws.Cells.Item(122, 9).Value = 11088  # I122
$ws.Cells.Item(122, 10).Value = 6124.875  # J122
$ws.Cells.Item(122, 11).Value = 33264  # K122
$ws.Cells.Item(122, 12).Value = 18374.625  # L122
$ws.Cells.Item(122, 8).Value = 8882.166999999999  # H122
$ws.Cells.Item(132, 12).Value = 10843.5  # L132
$ws.Cells.Item(132, 9).Value = 38462600  # I132
$ws.Cells.Item(132, 11).Value = 115387800  # K132
$ws.Cells.Item(132, 10).Value = 3614.5  # J132
$ws.Cells.Item(132, 8).Value = 33334736  # H132
$ws.Cells.Item(132, 13).Value = -115385270  # M132
$ws.Cells.Item(132, 14).Value = -15903.5  # N132
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(7, 8).Value = 5505.4  # H7
$ws.Cells.Item(7, 9).Value = 5175.6665  # I7
$ws.Cells.Item(7, 13).Value = -5063.6665  # M7
$ws.Cells.Item(7, 14).Value = -6224  # N7
$ws.Cells.Item(7, 12).Value = 6000  # L7
$ws.Cells.Item(7, 11).Value = 5175.6665  # K7
$ws.Cells.Item(7, 10).Value = 6000  # J7
$ws.Cells.Item(45, 8).Value = 0  # H45
$ws.Cells.Item(45, 11).Value = 0  # K45
$ws.Cells.Item(45, 13).Value = $null  # M45
$ws.Cells.Item(45, 9).Value = 0  # I45
$ws.Cells.Item(61, 13).Value = -2385.4167  # M61
$ws.Cells.Item(61, 12).Value = 5531.9287  # L61
$ws.Cells.Item(61, 14).Value = -5935.9287  # N61
$ws.Cells.Item(61, 10).Value = 5531.9287  # J61
$ws.Cells.Item(61, 11).Value = 2587.4167  # K61
$ws.Cells.Item(61, 8).Value = 4172.923  # H61
$ws.Cells.Item(61, 9).Value = 2587.4167  # I61
$ws.Cells.Item(62, 8).Value = 333359260  # H62
$ws.Cells.Item(62, 12).Value = 333359260  # L62
$ws.Cells.Item(62, 14).Value = -333360508  # N62
$ws.Cells.Item(62, 10).Value = 333359260  # J62
$ws.Cells.Item(65, 8).Value = 333359260  # H65
$ws.Cells.Item(65, 10).Value = 333359260  # J65
$ws.Cells.Item(65, 12).Value = 1000077780  # L65
$ws.Cells.Item(65, 14).Value = -1000084020  # N65
$ws.Cells.Item(93, 9).Value = 2835.6667  # I93
$ws.Cells.Item(93, 11).Value = 2835.6667  # K93
$ws.Cells.Item(93, 8).Value = 4702.4053  # H93
$ws.Cells.Item(93, 13).Value = -1587.6667  # M93
$ws.Cells.Item(113, 8).Value = 4172.923  # H113
$ws.Cells.Item(113, 13).Value = -417.4167000000002  # M113
$ws.Cells.Item(113, 11).Value = 2587.4167  # K113
$ws.Cells.Item(113, 9).Value = 2587.4167  # I113
$ws.Cells.Item(113, 12).Value = 5531.9287  # L113
$ws.Cells.Item(113, 10).Value = 5531.9287  # J113
$ws.Cells.Item(113, 14).Value = -9871.9287  # N113
$ws.Cells.Item(122, 14).Value = -22537.6  # N122
$ws.Cells.Item(122, 13).Value = -11195.875  # M122
$ws.Cells.Item(122, 9).Value = 4548.625  # I122
$ws.Cells.Item(122, 10).Value = 5879.2  # J122
$ws.Cells.Item(122, 11).Value = 13645.875  # K122
$ws.Cells.Item(122, 12).Value = 17637.6  # L122
$ws.Cells.Item(122, 8).Value = 5060.385  # H122
$ws.Cells.Item(126, 10).Value = 6000  # J126
$ws.Cells.Item(126, 9).Value = 5175.6665  # I126
$ws.Cells.Item(126, 14).Value = -22940  # N126
$ws.Cells.Item(126, 8).Value = 5505.4  # H126
$ws.Cells.Item(126, 11).Value = 15526.9995  # K126
$ws.Cells.Item(126, 13).Value = -13056.9995  # M126
$ws.Cells.Item(126, 12).Value = 18000  # L126
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(56, 10).Value = 59999  # J56
$ws.Cells.Item(56, 12).Value = 59999  # L56
$ws.Cells.Item(56, 13).Value = -13066  # M56
$ws.Cells.Item(56, 9).Value = 13780  # I56
$ws.Cells.Item(56, 8).Value = 36889.5  # H56
$ws.Cells.Item(56, 14).Value = -61427  # N56
$ws.Cells.Item(56, 11).Value = 13780  # K56
$ws.Cells.Item(81, 8).Value = 1268.2106  # H81
$ws.Cells.Item(81, 12).Value = 2162.5  # L81
$ws.Cells.Item(81, 10).Value = 1081.25  # J81
$ws.Cells.Item(81, 14).Value = -4284.5  # N81
$ws.Cells.Item(82, 8).Value = 504000000  # H82
$ws.Cells.Item(82, 11).Value = 0  # K82
$ws.Cells.Item(82, 9).Value = 0  # I82
$ws.Cells.Item(82, 13).Value = $null  # M82
$ws.Cells.Item(84, 12).Value = 10812.5  # L84
$ws.Cells.Item(84, 14).Value = -21420.5  # N84
$ws.Cells.Item(84, 10).Value = 1081.25  # J84
$ws.Cells.Item(84, 8).Value = 1268.2106  # H84
$ws.Cells.Item(85, 11).Value = 0  # K85
$ws.Cells.Item(85, 13).Value = $null  # M85
$ws.Cells.Item(85, 8).Value = 504000000  # H85
$ws.Cells.Item(85, 9).Value = 0  # I85
$ws.Cells.Item(88, 10).Value = 333366660  # J88
$ws.Cells.Item(88, 12).Value = 333366660  # L88
$ws.Cells.Item(88, 14).Value = -333367472  # N88
$ws.Cells.Item(88, 8).Value = 333366660  # H88
$ws.Cells.Item(91, 10).Value = 333366660  # J91
$ws.Cells.Item(91, 8).Value = 333366660  # H91
$ws.Cells.Item(91, 14).Value = -333369468  # N91
$ws.Cells.Item(91, 12).Value = 333366660  # L91
